# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Poroto verde, Femacal de La Calera) above
# the current row 303, shifting the existing data rows down, then populate
# the two new rows with the newly reported observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 303 (existing rows 303.. shift down to 305..)
$ws.Rows.Item(303).Resize(2).Insert()

# --- New row 303 ---
$ws.Cells.Item(303, 1).Value = 3
$ws.Cells.Item(303, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(303, 3).Value = "Coquimbo"
$ws.Cells.Item(303, 4).Value = 44736
$ws.Cells.Item(303, 5).Value = 5
$ws.Cells.Item(303, 6).Value = 100112031
$ws.Cells.Item(303, 7).Value = "Poroto verde"
$ws.Cells.Item(303, 8).Value = "Magnum"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 38
$ws.Cells.Item(303, 11).Value = 28000
$ws.Cells.Item(303, 12).Value = 28000
$ws.Cells.Item(303, 13).Value = 28000
$ws.Cells.Item(303, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(303, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(303, 16).Value = 1120
$ws.Cells.Item(303, 17).Value = 25
$ws.Cells.Item(303, 18).Value = "Hortaliza"

# --- New row 304 ---
$ws.Cells.Item(304, 1).Value = 3
$ws.Cells.Item(304, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(304, 3).Value = "Coquimbo"
$ws.Cells.Item(304, 4).Value = 44736
$ws.Cells.Item(304, 5).Value = 5
$ws.Cells.Item(304, 6).Value = 100112031
$ws.Cells.Item(304, 7).Value = "Poroto verde"
$ws.Cells.Item(304, 8).Value = "Magnum"
$ws.Cells.Item(304, 9).Value = "Primera"
$ws.Cells.Item(304, 10).Value = 40
$ws.Cells.Item(304, 11).Value = 23000
$ws.Cells.Item(304, 12).Value = 23000
$ws.Cells.Item(304, 13).Value = 23000
$ws.Cells.Item(304, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(304, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(304, 16).Value = 920
$ws.Cells.Item(304, 17).Value = 25
$ws.Cells.Item(304, 18).Value = "Hortaliza"
